$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.969.12"
$ws.Range("D3").Value = "2.417.68"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.18"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("E6").Value = "  -0.86%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.55%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  -0.86%  "
$ws.Range("E11").Value = "  -4.10%  "
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "2.854.01"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "61.903.79"
$ws.Range("E16").Value = "  -0.36%  "
$ws.Range("D17").Value = "2.452.60"
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "323.98"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.82"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.12"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.75"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.72"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.69"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "553.14"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.87%  "
$ws.Range("D27").Value = "2.534.58"
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("E31").Value = "  -4.90%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("E37").Value = "  -1.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "153.58"
$ws.Range("D38").ClearFormats()
$ws.Range("E39").Value = "  -5.86%  "
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("E44").Value = "  -6.52%  "
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0526"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.77%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.591"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("E51").Value = "  +4.37%  "
